$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (D2, J2) - K2/L2/M2/P2 remain unchanged
$ws.Range("D2").Value = 44291
$ws.Range("J2").Value = 30

# Update row 3 (D3, J3, K3, L3, M3, P3)
$ws.Range("D3").Value = 44277
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 11000
$ws.Range("P3").Value = 550

# Update row 5 (D5, J5, K5, L5, M5, P5)
$ws.Range("D5").Value = 44284
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 500
